# Daily attendance processing - 2025-12-06 15:24:43
# Reorders the list of "Recorded By" names in column G of each data row.
# The names (comma-separated) are sorted case-insensitively; when two
# names are equal ignoring case, the lexicographically "larger" one
# (e.g. lowercase variants) sorts first.

function CompareNames($a, $b) {
    $al = $a.ToLower()
    $bl = $b.ToLower()
    if ($al -ne $bl) {
        if ($al.CompareTo($bl) -lt 0) { return -1 } else { return 1 }
    }
    if ($a.CompareTo($b) -lt 0) { return 1 } else { return -1 }
}

function SortNames($arr) {
    for ($i = 1; $i -lt $arr.Count; $i++) {
        $key = $arr[$i]
        $j = $i - 1
        while ($j -ge 0 -and (CompareNames $arr[$j] $key) -gt 0) {
            $arr[$j + 1] = $arr[$j]
            $j = $j - 1
        }
        $arr[$j + 1] = $key
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -notlike "*,*") { continue }

    $rawParts = $val.Split(",")
    $names = @()
    foreach ($p in $rawParts) {
        $names += $p.Trim()
    }

    $sortedNames = SortNames $names
    $newVal = [string]::Join(", ", $sortedNames)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
